$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet has a stray leftover row (row 16: "Sheet", 3, 4)
# that is no longer needed. Select it and delete the entire row, which shifts every
# row below it up by one.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Activate()
$wsParams.Rows.Item(16).Select() | Out-Null
$wsParams.Rows.Item(16).Delete()

# Finish up by switching over to the "threshold_b" sheet/cell before saving.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select() | Out-Null
